$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the bug-description text in B4: append a bold "[FIXED]" tag ---
$cell = $ws.Range("B4")
$originalText = $cell.Value2
$suffix = "[FIXED]"
$newText = $originalText + " " + $suffix
$cell.Value2 = $newText

# Style the un-bolded portion (original sentence + trailing space) explicitly,
# matching the cell's existing Arial / black text font.
$bodyLen = $newText.Length - $suffix.Length
$bodyRange = $cell.Characters(1, $bodyLen)
$bodyRange.Font.Name = "Arial"
$bodyRange.Font.Color = 0

# Bold just the appended "[FIXED]" marker.
$fixedRange = $cell.Characters($bodyLen + 1, $suffix.Length)
$fixedRange.Font.Name = "Arial"
$fixedRange.Font.Color = 0
$fixedRange.Font.Bold = $true

# --- Enable word-wrap on the bug/expected columns for row 4 and grow the row ---
$ws.Range("B4:C4").WrapText = $true
$ws.Range("B4").EntireRow.RowHeight = 39
